$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 4035.3704
$ws.Range("I138").Value = 1478.2354
$ws.Range("J138").Value = 5885.213
$ws.Range("K138").Value = 4434.706200000001
$ws.Range("L138").Value = 17655.639
$ws.Range("M138").Value = 705.2937999999995
$ws.Range("N138").Value = -27935.639
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H112").Value = 25446.5
$ws.Range("J112").Value = 25446.5
$ws.Range("L112").Value = 25446.5
$ws.Range("N112").Value = -28400.5
$ws.Range("H119").Value = 40940.715
$ws.Range("J119").Value = 40940.715
$ws.Range("L119").Value = 40940.715
$ws.Range("N119").Value = -50616.715
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 64604
$ws.Range("I86").Value = 1070
$ws.Range("J86").Value = 128138
$ws.Range("K86").Value = 1070
$ws.Range("L86").Value = 128138
$ws.Range("M86").Value = 53
$ws.Range("N86").Value = -130384
$ws.Range("H89").Value = 64604
$ws.Range("I89").Value = 1070
$ws.Range("J89").Value = 128138
$ws.Range("K89").Value = 5350
$ws.Range("L89").Value = 640690
$ws.Range("M89").Value = 266
$ws.Range("N89").Value = -651922
$ws.Range("H133").Value = 20780
$ws.Range("J133").Value = 20780
$ws.Range("L133").Value = 20780
$ws.Range("N133").Value = -30900
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 30000
$ws.Range("J43").Value = 30000
$ws.Range("L43").Value = 30000
$ws.Range("N43").Value = -30368
$ws.Range("H58").Value = 18523568
$ws.Range("I58").Value = 3194.1538
$ws.Range("J58").Value = 35721056
$ws.Range("K58").Value = 3194.1538
$ws.Range("L58").Value = 35721056
$ws.Range("M58").Value = -2991.1538
$ws.Range("N58").Value = -35721462
$ws.Range("H92").Value = 40000
$ws.Range("J92").Value = 40000
$ws.Range("L92").Value = 40000
$ws.Range("N92").Value = -44992
$ws.Range("H99").Value = 4472.125
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 5296.1665
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 5296.1665
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -8292.166499999999
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H122").Value = 3011.4
$ws.Range("I122").Value = 2501.6
$ws.Range("J122").Value = 4540.8
$ws.Range("K122").Value = 7504.799999999999
$ws.Range("L122").Value = 13622.4
$ws.Range("M122").Value = -5054.799999999999
$ws.Range("N122").Value = -18522.4
$ws.Range("H126").Value = 4472.125
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 5296.1665
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 15888.4995
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -20828.4995
$ws.Range("H134").Value = 12197756
$ws.Range("I134").Value = 13515514
$ws.Range("J134").Value = 8499.75
$ws.Range("K134").Value = 40546542
$ws.Range("L134").Value = 25499.25
$ws.Range("M134").Value = -40544007
$ws.Range("N134").Value = -30569.25
$ws.Range("H136").Value = 18523568
$ws.Range("I136").Value = 3194.1538
$ws.Range("J136").Value = 35721056
$ws.Range("K136").Value = 9582.4614
$ws.Range("L136").Value = 107163168
$ws.Range("M136").Value = -7032.4614
$ws.Range("N136").Value = -107168268
$ws.Range("H139").Value = 20770
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("H141").Value = 28700
$ws.Range("J141").Value = 28700
$ws.Range("L141").Value = 28700
$ws.Range("N141").Value = -39060
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 263.33334
$ws.Range("I15").Value = 145
$ws.Range("K15").Value = 435
$ws.Range("M15").Value = -295
$ws.Range("H34").Value = 10360.818
$ws.Range("I34").Value = 160
$ws.Range("J34").Value = 16189.857
$ws.Range("K34").Value = 480
$ws.Range("L34").Value = 48569.571
$ws.Range("M34").Value = -396
$ws.Range("N34").Value = -48737.571
$ws.Range("H39").Value = 2447.6924
$ws.Range("J39").Value = 2447.6924
$ws.Range("L39").Value = 7343.0772
$ws.Range("N39").Value = -7931.0772
$ws.Range("H56").Value = 5462.3076
$ws.Range("I56").Value = 5462.3076
$ws.Range("K56").Value = 5462.3076
$ws.Range("M56").Value = -4932.3076
$ws.Range("H131").Value = 879.9036
$ws.Range("I131").Value = 401.81818
$ws.Range("J131").Value = 952.94446
$ws.Range("K131").Value = 1205.45454
$ws.Range("L131").Value = 2858.83338
$ws.Range("M131").Value = 3834.54546
$ws.Range("N131").Value = -12938.83338
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 30000
$ws.Range("J101").Value = 30000
$ws.Range("L101").Value = 30000
$ws.Range("N101").Value = -36490
$ws.Range("H102").Value = 49317.684
$ws.Range("I102").Value = 2940
$ws.Range("K102").Value = 2940
$ws.Range("M102").Value = -1318
$ws.Range("H122").Value = 5608.7026
$ws.Range("I122").Value = 5379.7646
$ws.Range("J122").Value = 5803.3
$ws.Range("K122").Value = 16139.2938
$ws.Range("L122").Value = 17409.9
$ws.Range("M122").Value = -13689.2938
$ws.Range("N122").Value = -22309.9
$ws.Range("H126").Value = 837602.2
$ws.Range("I126").Value = 2803
$ws.Range("J126").Value = 1255001.8
$ws.Range("K126").Value = 8409
$ws.Range("L126").Value = 3765005.4
$ws.Range("M126").Value = -5939
$ws.Range("N126").Value = -3769945.4
$ws.Range("H132").Value = 4693.4287
$ws.Range("I132").Value = 7102.6665
$ws.Range("J132").Value = 3552.2104
$ws.Range("K132").Value = 21307.9995
$ws.Range("L132").Value = 10656.6312
$ws.Range("M132").Value = -18777.9995
$ws.Range("N132").Value = -15716.6312
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2309.3157
$ws.Range("I7").Value = 1890.8572
$ws.Range("K7").Value = 1890.8572
$ws.Range("M7").Value = -1778.8572
$ws.Range("H40").Value = 5211.4116
$ws.Range("I40").Value = 5850.3335
$ws.Range("J40").Value = 3678
$ws.Range("K40").Value = 5850.3335
$ws.Range("L40").Value = 3678
$ws.Range("M40").Value = -5714.3335
$ws.Range("N40").Value = -3950
$ws.Range("H93").Value = 3610.3333
$ws.Range("I93").Value = 2440.5
$ws.Range("J93").Value = 5950
$ws.Range("K93").Value = 2440.5
$ws.Range("L93").Value = 5950
$ws.Range("M93").Value = -1192.5
$ws.Range("N93").Value = -8446
$ws.Range("H100").Value = 2450.923
$ws.Range("I100").Value = 1630.5
$ws.Range("J100").Value = 2815.5557
$ws.Range("K100").Value = 1630.5
$ws.Range("L100").Value = 2815.5557
$ws.Range("M100").Value = -1089.5
$ws.Range("N100").Value = -3897.5557
$ws.Range("H114").Value = 34000
$ws.Range("J114").Value = 34000
$ws.Range("L114").Value = 34000
$ws.Range("N114").Value = -42678
$ws.Range("H119").Value = 49800
$ws.Range("J119").Value = 49800
$ws.Range("L119").Value = 49800
$ws.Range("N119").Value = -59476
$ws.Range("H120").Value = 29733.334
$ws.Range("J120").Value = 29733.334
$ws.Range("L120").Value = 29733.334
$ws.Range("N120").Value = -39409.334
$ws.Range("H122").Value = 3160.9697
$ws.Range("I122").Value = 2690.1052
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 8070.3156
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -5620.3156
$ws.Range("N122").Value = -16300
$ws.Range("H126").Value = 2309.3157
$ws.Range("I126").Value = 1890.8572
$ws.Range("K126").Value = 5672.571599999999
$ws.Range("M126").Value = -3202.571599999999
$ws.Range("H134").Value = 47400
$ws.Range("J134").Value = 47400
$ws.Range("L134").Value = 47400
$ws.Range("N134").Value = -57540
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H111").Value = 28333.334
$ws.Range("J111").Value = 28333.334
$ws.Range("L111").Value = 28333.334
$ws.Range("N111").Value = -36513.334
$ws.Range("H122").Value = 1953.16
$ws.Range("I122").Value = 1491.1052
$ws.Range("J122").Value = 3416.3333
$ws.Range("K122").Value = 4473.3156
$ws.Range("L122").Value = 10248.9999
$ws.Range("M122").Value = -2023.3156
$ws.Range("N122").Value = -15148.9999
$ws.Range("H136").Value = 2312.95
$ws.Range("I136").Value = 1752.138
$ws.Range("K136").Value = 5256.414
$ws.Range("M136").Value = -2706.414
